$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "posx1 / posx2" mini-table mirroring the existing B30:C31 check (rows 22-23, cols H-I)
$ws.Range("H22").Value = "posx1"
$ws.Range("I22").Value = "posx2"
$ws.Range("H23").Value = 2.460673809
$ws.Range("I23").Value = -0.000000238

# "Distx" label and its computed difference
$ws.Range("H25").Value = "Distx"
$ws.Range("H26").Formula = "=H23-I23"

# New "Posy1 / Posy2" mini-table mirroring the existing B34:C35 check (rows 28-29, cols H-I)
$ws.Range("H28").Value = "Posy1"
$ws.Range("I28").Value = "Posy2"
$ws.Range("H29").Value = 0.000000229
$ws.Range("I29").Value = 4.262012482

# "Disty" computed difference
$ws.Range("H30").Formula = "=I29-H29"

# Clean up stray formatted-but-empty cell
$ws.Range("B33").ClearContents()

# Leave selection where the author finished reviewing
$ws.Range("J27").Select()
